# The "Skills Imperative 2035" release-date caveat/error has now been fixed
# upstream, so replace the previous (incorrect) placeholder date with the
# corrected release date for the "Latest period (release date)" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("C14").Value = "2035 (02/08/24)"

# Leave the view scrolled down a little with the selection left on B14,
# matching the saved workbook's recorded view state.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("B14").Select()
